$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.252.53"
$ws.Range("E2").Value = "  +4.52%  "
$ws.Range("D3").Value = "3.654.18"
$ws.Range("E3").Value = "  +9.65%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "644.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.48"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.403"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.66%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("E10").Value = "  +4.56%  "
$ws.Range("D11").Value = "3.649.05"
$ws.Range("E11").Value = "  +9.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("E13").Value = "  +3.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.45%  "
$ws.Range("D15").Value = "4.350.22"
$ws.Range("E15").Value = "  +9.93%  "
$ws.Range("D16").Value = "96.258.44"
$ws.Range("E16").Value = "  +4.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000257"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +24.94%  "
$ws.Range("D20").Value = "3.663.36"
$ws.Range("E20").Value = "  +9.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "519.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.00%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.484"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000198"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "98.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.49%  "
$ws.Range("E29").Value = "  +21.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.58%  "
$ws.Range("E31").Value = "  +2.86%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.181"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.11%  "
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "32.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +13.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.579"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "566.24"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.944"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.53%  "
$ws.Range("E41").Value = "  +2.94%  "
$ws.Range("B43").Value = "ImmutableX"
$ws.Range("C43").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.74%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0430"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.92%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("E47").Value = "  +5.67%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.67%  "
$ws.Range("B49").Value = "MantraDAO"
$ws.Range("C49").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.44%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "207.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +15.00%  "
